# Natmi following Dr Hou advice
# The LR-pairs sheet (Il34-Ptprz1) is recomputed to include a 3rd cluster,
# "ECs", alongside the existing "FAPs" and "sCs" clusters. The data block
# becomes the full 3x3 cross join of Sending cluster x Target cluster
# (6 rows instead of 4), and several numeric columns are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E..T).
$rows = @(
  @("FAPs","Il34","Ptprz1","sCs",
    2,0.6666666666666666,1.393664333333333,4.180993,0.04358216200908328,0.04358216200908328,
    3,1,0.1030763333333333,0.309229,0.01126512502660735,0.01126512502660735,
    0.1436538093774444,1.292884284397,0.0004909585039621803,0.0004909585039621803),
  @("FAPs","Il34","Ptprz1","ECs",
    2,0.6666666666666666,1.393664333333333,4.180993,0.04358216200908328,0.04358216200908328,
    3,1,9.046962666666667,27.140888,0.9887348749733926,0.9887348749733927,
    12.60842919353156,113.475862741784,0.0430912035051211,0.0430912035051211),
  @("sCs","Il34","Ptprz1","sCs",
    3,1,2.952820666666666,8.858461999999999,0.092339529397755,0.092339529397755,
    3,1,0.1030763333333333,0.309229,0.01126512502660735,0.01126512502660735,
    0.3043659273108889,2.739293345798,0.001040216343563795,0.001040216343563795),
  @("sCs","Il34","Ptprz1","ECs",
    3,1,2.952820666666666,8.858461999999999,0.092339529397755,0.092339529397755,
    3,1,9.046962666666667,27.140888,0.9887348749733926,0.9887348749733927,
    26.71405833269511,240.426524994256,0.0912993130541912,0.09129931305419121),
  @("ECs","Il34","Ptprz1","sCs",
    3,1,27.63137633333334,82.89412900000001,0.8640783085931617,0.8640783085931617,
    3,1,0.1030763333333333,0.309229,0.01126512502660735,0.01126512502660735,
    2.848140957393444,25.633268616541,0.009733950179081376,0.009733950179081378),
  @("ECs","Il34","Ptprz1","ECs",
    3,1,27.63137633333334,82.89412900000001,0.8640783085931617,0.8640783085931617,
    3,1,9.046962666666667,27.140888,0.9887348749733926,0.9887348749733927,
    249.9800301162836,2249.820271046552,0.8543443584140803,0.8543443584140804)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $line = $rows[$i]
  for ($c = 1; $c -le $line.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $line[$c - 1]
  }
}
